$p = $ppt.ActivePresentation
$d = $p.Designs.Item(1)
Write-Output $d.Name
$d.Name = "Office Theme"
$d2 = $p.Designs.Item(1)
Write-Output $d2.Name
